$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 94
$ws.Range("B94").Value = 6236252
$ws.Range("F94").Value = "Deportivo Tachira"
$ws.Range("G94").Value = "CD Hermanos Colmenares"
$ws.Range("H94").Value = 1
$ws.Range("J94").Value = "H"
$ws.Range("K94").Value = 1.363
$ws.Range("L94").Value = 4.2
$ws.Range("M94").Value = 7.5
$ws.Range("N94").Value = 1.333
$ws.Range("O94").Value = 4.5
$ws.Range("P94").Value = 8
$ws.Range("Q94").Value = -1.5
$ws.Range("R94").Value = 2
$ws.Range("S94").Value = 1.8
$ws.Range("T94").Value = 2.5
$ws.Range("U94").Value = 1.925
$ws.Range("V94").Value = 1.875
$ws.Range("W94").Value = 0.333
$ws.Range("X94").Value = -1
$ws.Range("Z94").Value = -1
$ws.Range("AA94").Value = 0.8
$ws.Range("AC94").Value = 0.875

# Row 95
$ws.Range("B95").Value = 6236255
$ws.Range("F95").Value = "Deportivo Rayo Zuliano"
$ws.Range("G95").Value = "Caracas"
$ws.Range("I95").Value = 0
$ws.Range("J95").Value = "D"
$ws.Range("K95").Value = 3.75
$ws.Range("M95").Value = 1.95
$ws.Range("N95").Value = 2.9
$ws.Range("O95").Value = 2.875
$ws.Range("P95").Value = 2.45
$ws.Range("Q95").Value = 0.25
$ws.Range("R95").Value = 1.775
$ws.Range("S95").Value = 2.025
$ws.Range("U95").Value = 1.85
$ws.Range("V95").Value = 1.95
$ws.Range("X95").Value = 1.875
$ws.Range("Y95").Value = -1
$ws.Range("Z95").Value = 0.3875
$ws.Range("AA95").Value = -0.5
$ws.Range("AB95").Value = -1
$ws.Range("AC95").Value = 0.95

# Row 96
$ws.Range("B96").Value = 6236251
$ws.Range("F96").Value = "Angostura FC"
$ws.Range("G96").Value = "Portuguesa"
$ws.Range("I96").Value = 2
$ws.Range("J96").Value = "A"
$ws.Range("K96").Value = 3.1
$ws.Range("L96").Value = 3.2
$ws.Range("M96").Value = 2.15
$ws.Range("N96").Value = 4
$ws.Range("O96").Value = 3.6
$ws.Range("P96").Value = 1.75
$ws.Range("Q96").Value = 0.75
$ws.Range("R96").Value = 1.8
$ws.Range("S96").Value = 2
$ws.Range("U96").Value = 1.95
$ws.Range("V96").Value = 1.85
$ws.Range("W96").Value = -1
$ws.Range("Y96").Value = 0.75
$ws.Range("Z96").Value = -0.5
$ws.Range("AA96").Value = 0.5
$ws.Range("AB96").Value = 0.95
$ws.Range("AC96").Value = -1

# Row 97
$ws.Range("B97").Value = 6236253
$ws.Range("F97").Value = "Deportivo La Guaira"
$ws.Range("G97").Value = "UCV"
$ws.Range("H97").Value = 0
$ws.Range("J97").Value = "D"
$ws.Range("K97").Value = 1.833
$ws.Range("L97").Value = 3.25
$ws.Range("M97").Value = 4
$ws.Range("N97").Value = 2
$ws.Range("O97").Value = 3.2
$ws.Range("P97").Value = 3.5
$ws.Range("Q97").Value = -0.25
$ws.Range("R97").Value = 1.775
$ws.Range("S97").Value = 2.025
$ws.Range("T97").Value = 2.25
$ws.Range("U97").Value = 1.9
$ws.Range("V97").Value = 1.9
$ws.Range("W97").Value = -1
$ws.Range("X97").Value = 2.2
$ws.Range("Z97").Value = -0.5
$ws.Range("AA97").Value = 0.5125
$ws.Range("AC97").Value = 0.8999999999999999

# Row 98
$ws.Range("B98").Value = 6236254
$ws.Range("F98").Value = "Academia Puerto Cabello"
$ws.Range("G98").Value = "Estudiantes Merida"
$ws.Range("I98").Value = 0
$ws.Range("J98").Value = "H"
$ws.Range("K98").Value = 1.727
$ws.Range("L98").Value = 3.4
$ws.Range("M98").Value = 4.333
$ws.Range("N98").Value = 1.666
$ws.Range("O98").Value = 3.4
$ws.Range("P98").Value = 4.75
$ws.Range("Q98").Value = -0.75
$ws.Range("R98").Value = 1.875
$ws.Range("S98").Value = 1.925
$ws.Range("U98").Value = 1.9
$ws.Range("V98").Value = 1.9
$ws.Range("W98").Value = 0.6659999999999999
$ws.Range("Y98").Value = -1
$ws.Range("Z98").Value = 0.4375
$ws.Range("AA98").Value = -0.5
$ws.Range("AB98").Value = -1
$ws.Range("AC98").Value = 0.8999999999999999

# Row 99
$ws.Range("B99").Value = 6236612
$ws.Range("F99").Value = "Zamora"
$ws.Range("G99").Value = "Carabobo"
$ws.Range("I99").Value = 2
$ws.Range("J99").Value = "A"
$ws.Range("K99").Value = 3.2
$ws.Range("M99").Value = 2.15
$ws.Range("N99").Value = 4.5
$ws.Range("O99").Value = 3.3
$ws.Range("P99").Value = 1.75
$ws.Range("Q99").Value = 0.5
$ws.Range("R99").Value = 2
$ws.Range("S99").Value = 1.8
$ws.Range("U99").Value = 1.925
$ws.Range("V99").Value = 1.875
$ws.Range("X99").Value = -1
$ws.Range("Y99").Value = 0.75
$ws.Range("Z99").Value = -1
$ws.Range("AA99").Value = 0.8
$ws.Range("AB99").Value = -0.5
$ws.Range("AC99").Value = 0.4375

# Row 114
$ws.Range("B114").Value = 7352251
$ws.Range("F114").Value = "Caracas"
$ws.Range("G114").Value = "Academia Puerto Cabello"
$ws.Range("I114").Value = 0
$ws.Range("J114").Value = "H"
$ws.Range("K114").Value = 2.1
$ws.Range("L114").Value = 3.2
$ws.Range("M114").Value = 3.3
$ws.Range("N114").Value = 2.15
$ws.Range("O114").Value = 3.1
$ws.Range("P114").Value = 3.2
$ws.Range("Q114").Value = -0.5
$ws.Range("R114").Value = 2.025
$ws.Range("S114").Value = 1.775
$ws.Range("T114").Value = 2.25
$ws.Range("U114").Value = 1.975
$ws.Range("V114").Value = 1.825
$ws.Range("W114").Value = 1.15
$ws.Range("X114").Value = -1
$ws.Range("Z114").Value = 1.025
$ws.Range("AA114").Value = -1
$ws.Range("AB114").Value = -1
$ws.Range("AC114").Value = 0.825

# Row 115
$ws.Range("B115").Value = 7352250
$ws.Range("F115").Value = "Portuguesa"
$ws.Range("G115").Value = "Deportivo Tachira"
$ws.Range("I115").Value = 1
$ws.Range("J115").Value = "D"
$ws.Range("K115").Value = 3.1
$ws.Range("L115").Value = 2.875
$ws.Range("M115").Value = 2.3
$ws.Range("N115").Value = 3
$ws.Range("O115").Value = 2.875
$ws.Range("P115").Value = 2.375
$ws.Range("Q115").Value = 0.25
$ws.Range("R115").Value = 1.725
$ws.Range("S115").Value = 2.075
$ws.Range("T115").Value = 2
$ws.Range("U115").Value = 1.825
$ws.Range("V115").Value = 1.975
$ws.Range("W115").Value = -1
$ws.Range("X115").Value = 1.875
$ws.Range("Z115").Value = 0.3625
$ws.Range("AA115").Value = -0.5
$ws.Range("AB115").Value = 0
$ws.Range("AC115").Value = -0

# Row 116
$ws.Range("B116").Value = 7352252
$ws.Range("F116").Value = "Deportivo Tachira"
$ws.Range("G116").Value = "Caracas"
$ws.Range("K116").Value = 2.3
$ws.Range("L116").Value = 2.875
$ws.Range("M116").Value = 3.1
$ws.Range("N116").Value = 2.25
$ws.Range("O116").Value = 2.8
$ws.Range("P116").Value = 3.25
$ws.Range("R116").Value = 1.975
$ws.Range("S116").Value = 1.825
$ws.Range("T116").Value = 2
$ws.Range("U116").Value = 1.925
$ws.Range("V116").Value = 1.875
$ws.Range("X116").Value = 1.8
$ws.Range("AA116").Value = 0.4125
$ws.Range("AB116").Value = 0
$ws.Range("AC116").Value = -0

# Row 117
$ws.Range("B117").Value = 7352254
$ws.Range("F117").Value = "Academia Puerto Cabello"
$ws.Range("G117").Value = "Portuguesa"
$ws.Range("K117").Value = 2.05
$ws.Range("L117").Value = 3.4
$ws.Range("M117").Value = 3
$ws.Range("N117").Value = 1.833
$ws.Range("O117").Value = 3.5
$ws.Range("P117").Value = 3.5
$ws.Range("R117").Value = 1.65
$ws.Range("S117").Value = 2.2
$ws.Range("T117").Value = 2.25
$ws.Range("U117").Value = 1.825
$ws.Range("V117").Value = 1.975
$ws.Range("X117").Value = 2.5
$ws.Range("AA117").Value = 0.6000000000000001
$ws.Range("AB117").Value = -0.5
$ws.Range("AC117").Value = 0.4875

# Row 156
$ws.Range("H156").Value = 0
$ws.Range("I156").Value = 2
$ws.Range("J156").Value = "A"
$ws.Range("N156").Value = 2.1
$ws.Range("O156").Value = 3.2
$ws.Range("P156").Value = 3.25
$ws.Range("Q156").Value = -0.25
$ws.Range("R156").Value = 1.85
$ws.Range("S156").Value = 1.95
$ws.Range("U156").Value = 1.9
$ws.Range("V156").Value = 1.9
$ws.Range("W156").Value = -1
$ws.Range("X156").Value = -1
$ws.Range("Y156").Value = 2.25
$ws.Range("Z156").Value = -1
$ws.Range("AA156").Value = 0.95
$ws.Range("AB156").Value = -0.5
$ws.Range("AC156").Value = 0.45

# Row 157
$ws.Range("B157").Value = 7920997
$ws.Range("F157").Value = "Carabobo"
$ws.Range("G157").Value = "UCV"
$ws.Range("H157").Value = 0
$ws.Range("I157").Value = 1
$ws.Range("J157").Value = "A"
$ws.Range("K157").Value = 1.833
$ws.Range("L157").Value = 3.1
$ws.Range("M157").Value = 4.2
$ws.Range("N157").Value = 1.833
$ws.Range("O157").Value = 3.1
$ws.Range("P157").Value = 4.2
$ws.Range("Q157").Value = -0.5
$ws.Range("R157").Value = 1.9
$ws.Range("S157").Value = 1.9
$ws.Range("U157").Value = 1.85
$ws.Range("V157").Value = 1.95
$ws.Range("W157").Value = -1
$ws.Range("X157").Value = -1
$ws.Range("Y157").Value = 3.2
$ws.Range("Z157").Value = -1
$ws.Range("AA157").Value = 0.8999999999999999
$ws.Range("AB157").Value = -1
$ws.Range("AC157").Value = 0.95

# Row 158
$ws.Range("B158").Value = 7920998
$ws.Range("F158").Value = "Zamora"
$ws.Range("G158").Value = "Caracas"
$ws.Range("H158").Value = 2
$ws.Range("I158").Value = 2
$ws.Range("J158").Value = "D"
$ws.Range("K158").Value = 3.75
$ws.Range("L158").Value = 3.2
$ws.Range("M158").Value = 1.909
$ws.Range("N158").Value = 3
$ws.Range("O158").Value = 2.9
$ws.Range("P158").Value = 2.375
$ws.Range("Q158").Value = 0.25
$ws.Range("R158").Value = 1.8
$ws.Range("S158").Value = 2
$ws.Range("U158").Value = 1.825
$ws.Range("V158").Value = 1.975
$ws.Range("W158").Value = -1
$ws.Range("X158").Value = 1.9
$ws.Range("Y158").Value = -1
$ws.Range("Z158").Value = 0.4
$ws.Range("AA158").Value = -0.5
$ws.Range("AB158").Value = 0.825
$ws.Range("AC158").Value = -1

# Row 159
$ws.Range("H159").Value = 2
$ws.Range("I159").Value = 3
$ws.Range("J159").Value = "A"
$ws.Range("N159").Value = 3.4
$ws.Range("P159").Value = 2.2
$ws.Range("Q159").Value = 0.25
$ws.Range("R159").Value = 1.875
$ws.Range("S159").Value = 1.925
$ws.Range("U159").Value = 2.025
$ws.Range("V159").Value = 1.775
$ws.Range("W159").Value = -1
$ws.Range("X159").Value = -1
$ws.Range("Y159").Value = 1.2
$ws.Range("Z159").Value = -1
$ws.Range("AA159").Value = 0.925
$ws.Range("AB159").Value = 1.025
$ws.Range("AC159").Value = -1

# Row 160
$ws.Range("H160").Value = 1
$ws.Range("I160").Value = 2
$ws.Range("J160").Value = "A"
$ws.Range("N160").Value = 2.375
$ws.Range("P160").Value = 2.875
$ws.Range("Q160").Value = 0
$ws.Range("R160").Value = 1.75
$ws.Range("S160").Value = 2.05
$ws.Range("U160").Value = 1.9
$ws.Range("V160").Value = 1.9
$ws.Range("W160").Value = -1
$ws.Range("X160").Value = -1
$ws.Range("Y160").Value = 1.875
$ws.Range("Z160").Value = -1
$ws.Range("AA160").Value = 1.05
$ws.Range("AB160").Value = 0.8999999999999999
$ws.Range("AC160").Value = -1
